# feat: LootTable and CharacterStats
# Add a new "Items" table definition row (row 7) to the Luban _Tables sheet,
# mirroring the existing rows (Constants, Characters, Maps), and move the
# active selection to the newly added output-filename cell (E7).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New row 7: full_name / value_type / read_schema_from_file / output
$ws.Range("B7").Value = "Items.TbItems"
$ws.Range("C7").Value = "ItemTemplate"
$ws.Range("D7").Value = $true
$ws.Range("E7").Value = "Items.xlsx"

# Update the selected cell to the new output-filename cell, as in the diff.
$ws.Range("E7").Select() | Out-Null
